# Insert 3 new rows before row 18 (pushing the existing rows 18-23 down to 21-26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18:T20").EntireRow.Insert()

# Row 18: new "Especial" entry for the week of 44455
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44455
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = "Otros"
$ws.Range("I18").Value = 100107002
$ws.Range("J18").Value = "Chirimoya"
$ws.Range("K18").Value = "Cultivar IV Región"
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 340
$ws.Range("N18").Value = 2300
$ws.Range("O18").Value = 2400
$ws.Range("P18").Value = 2350
$ws.Range("Q18").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 2350
$ws.Range("T18").Value = 1

# Row 19: new "Extra (doble especial)" entry for the week of 44455
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44455
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = "Otros"
$ws.Range("I19").Value = 100107002
$ws.Range("J19").Value = "Chirimoya"
$ws.Range("K19").Value = "Cultivar IV Región"
$ws.Range("L19").Value = "Extra (doble especial)"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 2700
$ws.Range("O19").Value = 2800
$ws.Range("P19").Value = 2750
$ws.Range("Q19").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 2750
$ws.Range("T19").Value = 1

# Row 20: new "Primera" entry for the week of 44455
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44455
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107002
$ws.Range("J20").Value = "Chirimoya"
$ws.Range("K20").Value = "Cultivar IV Región"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 2000
$ws.Range("O20").Value = 2100
$ws.Range("P20").Value = 2050
$ws.Range("Q20").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 2050
$ws.Range("T20").Value = 1

# Update dimension-affecting used range is handled automatically by Excel;
# rows 21-26 already contain the shifted former rows 18-23 with correct values.
